$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8:10 (the "Desktop Computer" / DKS entries), shifting cells up
$ws.Range("A8:G10").EntireRow.Delete() | Out-Null

# Select cell E10, matching the final selection state
$ws.Range("E10").Select() | Out-Null

# Apply page setup (paper size A4, portrait orientation)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
